$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "clockify" right after DeskTime
$ws.Rows.Item(7).Insert() | Out-Null

# clockify row data
$ws.Range("A7").Value2 = "clockify"
$ws.Range("B7").Value2 = "✓"
$ws.Range("C7").Value2 = "✓"
$ws.Range("D7").Value2 = "!"
$ws.Range("E7").Value2 = "✓"
$ws.Range("F7").Value2 = "✓"
$ws.Range("G7").Value2 = "!"

# Row height / formatting consistent with the rest of the table
$ws.Rows.Item(7).RowHeight = 33.75
$ws.Range("A7:G7").HorizontalAlignment = -4108
$ws.Range("A7:G7").VerticalAlignment = -4108
$ws.Range("A7:G7").WrapText = $true

# Tool name is bold for the new row
$ws.Range("A7").Font.Bold = $true

# --- Header row: "Admite Release+Report+Burndown" -> "Admite Release+Report detailed"
$ws.Range("E1").Value2 = "Admite Release+Report detailed"

# --- Harvest row (row 4): "Gratis" column goes from "!" to "✗"
$ws.Range("D4").Value2 = "✗"

# --- Remove the DeskTime hyperlink (row 6)
$ws.Hyperlinks.Delete()

$ws.Range("E9").Select() | Out-Null
